$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update lasAHL (nM) values in column B (rows 2-13)
$ws.Range("B2").Value = 1621
$ws.Range("B3").Value = 1615
$ws.Range("B4").Value = 1550
$ws.Range("B5").Value = 1305
$ws.Range("B6").Value = 1257
$ws.Range("B7").Value = 575
$ws.Range("B8").Value = 254
$ws.Range("B9").Value = 45.9
$ws.Range("B10").Value = 22.7
$ws.Range("B11").Value = 5.64
$ws.Range("B12").Value = 4.6399999999999997
$ws.Range("B13").Value = 2.66

# Clear row 14 (A14 and B14) but keep formatting/style
$ws.Range("A14:B14").ClearContents()

# Update the selection to F10
$ws.Range("F10").Select()
